$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "test data 1"
$ws.Range("B3").Value = "test data 2"
$ws.Range("B4").Value = "test data 3"
$ws.Range("B5").Value = "test data 4"
$ws.Range("B6").Value = "test data 5"
$ws.Range("B7").Value = "test data 6"
$ws.Range("B8").Value = "test data 7"
$ws.Range("B9").Value = "test data 8"
$ws.Range("B10").Value = "test data 9"
$ws.Range("B11").Value = "test data 10"

$ws.Range("B13").Select()
